$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-09-01 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-02 Saturday", 2) | Out-Null

# Update the 100 arithmetic-problem table cells (row-major, 20 rows x 5 cols)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "42-0="
$t.Cell(1, 2).Range.Text = "77-12="
$t.Cell(1, 3).Range.Text = "75-13="
$t.Cell(1, 4).Range.Text = "53+9="
$t.Cell(1, 5).Range.Text = "57-2="
$t.Cell(2, 1).Range.Text = "65-53="
$t.Cell(2, 2).Range.Text = "86-80="
$t.Cell(2, 3).Range.Text = "1+3="
$t.Cell(2, 4).Range.Text = "25+22="
$t.Cell(2, 5).Range.Text = "57+36="
$t.Cell(3, 1).Range.Text = "79-71="
$t.Cell(3, 2).Range.Text = "37+33="
$t.Cell(3, 3).Range.Text = "7+63="
$t.Cell(3, 4).Range.Text = "50-5="
$t.Cell(3, 5).Range.Text = "95-69="
$t.Cell(4, 1).Range.Text = "32+15="
$t.Cell(4, 2).Range.Text = "83-13="
$t.Cell(4, 3).Range.Text = "73-5="
$t.Cell(4, 4).Range.Text = "34+59="
$t.Cell(4, 5).Range.Text = "94-21="
$t.Cell(5, 1).Range.Text = "8+3="
$t.Cell(5, 2).Range.Text = "17+78="
$t.Cell(5, 3).Range.Text = "24-4="
$t.Cell(5, 4).Range.Text = "44+55="
$t.Cell(5, 5).Range.Text = "18-9="
$t.Cell(6, 1).Range.Text = "18+7="
$t.Cell(6, 2).Range.Text = "99-87="
$t.Cell(6, 3).Range.Text = "7+9="
$t.Cell(6, 4).Range.Text = "39-26="
$t.Cell(6, 5).Range.Text = "20-4="
$t.Cell(7, 1).Range.Text = "18+52="
$t.Cell(7, 2).Range.Text = "47-40="
$t.Cell(7, 3).Range.Text = "70+14="
$t.Cell(7, 4).Range.Text = "21+24="
$t.Cell(7, 5).Range.Text = "22+17="
$t.Cell(8, 1).Range.Text = "36+40="
$t.Cell(8, 2).Range.Text = "15+82="
$t.Cell(8, 3).Range.Text = "69-39="
$t.Cell(8, 4).Range.Text = "37-18="
$t.Cell(8, 5).Range.Text = "34+26="
$t.Cell(9, 1).Range.Text = "96+2="
$t.Cell(9, 2).Range.Text = "89-25="
$t.Cell(9, 3).Range.Text = "23+2="
$t.Cell(9, 4).Range.Text = "54-52="
$t.Cell(9, 5).Range.Text = "75-25="
$t.Cell(10, 1).Range.Text = "89-27="
$t.Cell(10, 2).Range.Text = "88-29="
$t.Cell(10, 3).Range.Text = "32-23="
$t.Cell(10, 4).Range.Text = "37-0="
$t.Cell(10, 5).Range.Text = "98-44="
$t.Cell(11, 1).Range.Text = "26+12="
$t.Cell(11, 2).Range.Text = "73+22="
$t.Cell(11, 3).Range.Text = "37+28="
$t.Cell(11, 4).Range.Text = "56+29="
$t.Cell(11, 5).Range.Text = "73-23="
$t.Cell(12, 1).Range.Text = "75-23="
$t.Cell(12, 2).Range.Text = "16-15="
$t.Cell(12, 3).Range.Text = "55-15="
$t.Cell(12, 4).Range.Text = "5+7="
$t.Cell(12, 5).Range.Text = "26+40="
$t.Cell(13, 1).Range.Text = "93-67="
$t.Cell(13, 2).Range.Text = "88-54="
$t.Cell(13, 3).Range.Text = "35+17="
$t.Cell(13, 4).Range.Text = "43-30="
$t.Cell(13, 5).Range.Text = "30-25="
$t.Cell(14, 1).Range.Text = "99-75="
$t.Cell(14, 2).Range.Text = "15-1="
$t.Cell(14, 3).Range.Text = "84-13="
$t.Cell(14, 4).Range.Text = "61-47="
$t.Cell(14, 5).Range.Text = "84-11="
$t.Cell(15, 1).Range.Text = "75+7="
$t.Cell(15, 2).Range.Text = "95-16="
$t.Cell(15, 3).Range.Text = "33-3="
$t.Cell(15, 4).Range.Text = "69-11="
$t.Cell(15, 5).Range.Text = "33+30="
$t.Cell(16, 1).Range.Text = "39-3="
$t.Cell(16, 2).Range.Text = "0+21="
$t.Cell(16, 3).Range.Text = "0+72="
$t.Cell(16, 4).Range.Text = "36+13="
$t.Cell(16, 5).Range.Text = "71-32="
$t.Cell(17, 1).Range.Text = "17-1="
$t.Cell(17, 2).Range.Text = "1+55="
$t.Cell(17, 3).Range.Text = "98-72="
$t.Cell(17, 4).Range.Text = "33+25="
$t.Cell(17, 5).Range.Text = "48+46="
$t.Cell(18, 1).Range.Text = "81-24="
$t.Cell(18, 2).Range.Text = "99-22="
$t.Cell(18, 3).Range.Text = "38-4="
$t.Cell(18, 4).Range.Text = "38-31="
$t.Cell(18, 5).Range.Text = "60+1="
$t.Cell(19, 1).Range.Text = "64-19="
$t.Cell(19, 2).Range.Text = "28+8="
$t.Cell(19, 3).Range.Text = "18+47="
$t.Cell(19, 4).Range.Text = "23-1="
$t.Cell(19, 5).Range.Text = "96-80="
$t.Cell(20, 1).Range.Text = "14+33="
$t.Cell(20, 2).Range.Text = "4+71="
$t.Cell(20, 3).Range.Text = "14+36="
$t.Cell(20, 4).Range.Text = "17+59="
$t.Cell(20, 5).Range.Text = "44+0="
